# Actualización automática 2025-10-08 15:30:10
$wb = $excel.ActiveWorkbook

# --- Sheet 1: "VENTAS POR GRUPO" ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

$ws1.Range("H8").Value = 307.8
$ws1.Range("I8").Value = 129.6
$ws1.Range("M8").Value = 3555.5

$ws1.Range("H10").Value = 615.6
$ws1.Range("M10").Value = 1886.28

$ws1.Range("H15").Value = "3 de 13"
$ws1.Range("I15").Value = "1 de 13"
$ws1.Range("M15").Value = "2 de 13"

# --- Sheet 2: "VENTA MENSUAL" ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

$ws2.Range("F8").Value = 3992.9
$ws2.Range("F10").Value = 2501.88
$ws2.Range("F15").Value = 6512.77
